$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2098.6365
$ws.Range("J38").Value = 4960.3335
$ws.Range("L38").Value = 14881.0005
$ws.Range("N38").Value = -15625.0005
$ws.Range("H86").Value = 8685.733
$ws.Range("I86").Value = 12000
$ws.Range("J86").Value = 8449
$ws.Range("K86").Value = 12000
$ws.Range("L86").Value = 8449
$ws.Range("M86").Value = -10877
$ws.Range("N86").Value = -10695
$ws.Range("H89").Value = 8685.733
$ws.Range("I89").Value = 12000
$ws.Range("J89").Value = 8449
$ws.Range("K89").Value = 60000
$ws.Range("L89").Value = 42245
$ws.Range("M89").Value = -54384
$ws.Range("N89").Value = -53477
$ws.Range("H107").Value = 561.3125
$ws.Range("I107").Value = 587
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 587
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1333
$ws.Range("N107").Value = -4290
$ws.Range("H112").Value = 4279.4
$ws.Range("J112").Value = 5732.6665
$ws.Range("L112").Value = 17197.9995
$ws.Range("N112").Value = -19413.9995
$ws.Range("H138").Value = 4753.5854
$ws.Range("J138").Value = 5819.8
$ws.Range("L138").Value = 17459.4
$ws.Range("N138").Value = -27739.4

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3566.4666
$ws.Range("I122").Value = 3191.7693
$ws.Range("K122").Value = 9575.3079
$ws.Range("M122").Value = -7125.3079

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 5805.8335
$ws.Range("I35").Value = 709
$ws.Range("J35").Value = 15999.5
$ws.Range("K35").Value = 709
$ws.Range("L35").Value = 15999.5
$ws.Range("M35").Value = -415
$ws.Range("N35").Value = -16587.5
$ws.Range("H86").Value = 7211.5
$ws.Range("I86").Value = 6338.6
$ws.Range("K86").Value = 6338.6
$ws.Range("M86").Value = -5215.6
$ws.Range("H89").Value = 7211.5
$ws.Range("I89").Value = 6338.6
$ws.Range("K89").Value = 31693
$ws.Range("M89").Value = -26077
$ws.Range("H93").Value = 152044.33
$ws.Range("J93").Value = 167424.88
$ws.Range("L93").Value = 167424.88
$ws.Range("N93").Value = -171168.88
$ws.Range("H105").Value = 50001630
$ws.Range("J105").Value = 333335520
$ws.Range("L105").Value = 333335520
$ws.Range("N105").Value = -333339014
$ws.Range("H107").Value = 570.2778
$ws.Range("I107").Value = 452.35715
$ws.Range("K107").Value = 452.35715
$ws.Range("M107").Value = 1467.64285

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 625
$ws.Range("J5").Value = 400
$ws.Range("L5").Value = 1200
$ws.Range("N5").Value = -1424
$ws.Range("H14").Value = 1411
$ws.Range("I14").Value = 1411
$ws.Range("K14").Value = 4233
$ws.Range("M14").Value = -4060
$ws.Range("H38").Value = 102.55556
$ws.Range("I38").Value = 74.57143000000001
$ws.Range("J38").Value = 200.5
$ws.Range("K38").Value = 223.71429
$ws.Range("L38").Value = 601.5
$ws.Range("M38").Value = 123.28571
$ws.Range("N38").Value = -1295.5
$ws.Range("H68").Value = 4291.7144
$ws.Range("I68").Value = 5333.3335
$ws.Range("J68").Value = 4118.1113
$ws.Range("K68").Value = 16000.0005
$ws.Range("L68").Value = 12354.3339
$ws.Range("M68").Value = -15189.0005
$ws.Range("N68").Value = -13976.3339
$ws.Range("H71").Value = 4291.7144
$ws.Range("I71").Value = 5333.3335
$ws.Range("J71").Value = 4118.1113
$ws.Range("K71").Value = 48000.0015
$ws.Range("L71").Value = 37063.00169999999
$ws.Range("M71").Value = -43944.0015
$ws.Range("N71").Value = -45175.00169999999
$ws.Range("H82").Value = 17112.555
$ws.Range("I82").Value = 7006.5
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 21019.5
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -20613.5
$ws.Range("N82").Value = -60812
$ws.Range("H85").Value = 17112.555
$ws.Range("I85").Value = 7006.5
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 21019.5
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -19615.5
$ws.Range("N85").Value = -62808
$ws.Range("H113").Value = 901.35486
$ws.Range("J113").Value = 1021.36
$ws.Range("L113").Value = 3064.08
$ws.Range("N113").Value = -7404.08
$ws.Range("H131").Value = 758676.8
$ws.Range("I131").Value = 695.625
$ws.Range("K131").Value = 2086.875
$ws.Range("M131").Value = 2953.125
$ws.Range("H135").Value = 625
$ws.Range("J135").Value = 400
$ws.Range("L135").Value = 3600
$ws.Range("N135").Value = -8670

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3966.6038
$ws.Range("I70").Value = 3546.625
$ws.Range("J70").Value = 4041.2666
$ws.Range("K70").Value = 3546.625
$ws.Range("L70").Value = 4041.2666
$ws.Range("M70").Value = -3276.625
$ws.Range("N70").Value = -4581.2666
$ws.Range("H73").Value = 3966.6038
$ws.Range("I73").Value = 3546.625
$ws.Range("J73").Value = 4041.2666
$ws.Range("K73").Value = 3546.625
$ws.Range("L73").Value = 4041.2666
$ws.Range("M73").Value = -2610.625
$ws.Range("N73").Value = -5913.2666
$ws.Range("H107").Value = 523.5
$ws.Range("I107").Value = 398
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 398
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1522
$ws.Range("N107").Value = -4740
$ws.Range("H126").Value = 4801.1
$ws.Range("I126").Value = 3702.2
$ws.Range("K126").Value = 11106.6
$ws.Range("M126").Value = -8636.599999999999

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9155.5
$ws.Range("I7").Value = 10194.75
$ws.Range("K7").Value = 10194.75
$ws.Range("M7").Value = -10082.75
$ws.Range("H61").Value = 1795.027
$ws.Range("J61").Value = 2905.923
$ws.Range("L61").Value = 2905.923
$ws.Range("N61").Value = -3309.923
$ws.Range("H93").Value = 2265
$ws.Range("I93").Value = 2177.8333
$ws.Range("K93").Value = 2177.8333
$ws.Range("M93").Value = -929.8332999999998
$ws.Range("H113").Value = 1795.027
$ws.Range("J113").Value = 2905.923
$ws.Range("L113").Value = 2905.923
$ws.Range("N113").Value = -7245.923
$ws.Range("H126").Value = 9155.5
$ws.Range("I126").Value = 10194.75
$ws.Range("K126").Value = 30584.25
$ws.Range("M126").Value = -28114.25
$ws.Range("H132").Value = 4139.225
$ws.Range("I132").Value = 3718.9333
$ws.Range("J132").Value = 5400.1
$ws.Range("K132").Value = 11156.7999
$ws.Range("L132").Value = 16200.3
$ws.Range("M132").Value = -8626.7999
$ws.Range("N132").Value = -21260.3

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7599.8
$ws.Range("I81").Value = 1999.75
$ws.Range("J81").Value = 30000
$ws.Range("K81").Value = 3999.5
$ws.Range("L81").Value = 60000
$ws.Range("M81").Value = -2938.5
$ws.Range("N81").Value = -62122
$ws.Range("H84").Value = 7599.8
$ws.Range("I84").Value = 1999.75
$ws.Range("J84").Value = 30000
$ws.Range("K84").Value = 19997.5
$ws.Range("L84").Value = 300000
$ws.Range("M84").Value = -14693.5
$ws.Range("N84").Value = -310608
$ws.Range("H126").Value = 1765.6842
$ws.Range("I126").Value = 1610.0555
$ws.Range("K126").Value = 4830.166499999999
$ws.Range("M126").Value = -2360.166499999999
$ws.Range("H132").Value = 3396.8
$ws.Range("I132").Value = 2565.375
$ws.Range("K132").Value = 7696.125
$ws.Range("M132").Value = -5166.125
